# "Generate Report for Handback"
#
# The localization handback has completed: the Overview/zh-cn/de-de status
# columns flip from "Ready for handoff" to "Handed back: in sync with en-US",
# the per-language tables pick up the generated target + handback file names
# (with hyperlinks to the source .md on GitHub) and a fresh handback
# timestamp, and the "Latest Target/Handback File" columns are widened so the
# new file names are readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status cells -----------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$TargetFile2,
        [string]$TargetFile3,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File / Latest Handback File text (columns I / J).
    $ws.Range("I2").Value = "b9a358c7-8e1c-46ea-94ea-a5e869538c3e.md"
    $ws.Range("J2").Value = $TargetFile2

    $ws.Range("I3").Value = "e0a40405-df7f-4f03-93b3-560f8c411baf.md"
    $ws.Range("J3").Value = $TargetFile3

    # Latest Handback DateTime (column K).
    $ws.Range("K2").Value = $HandbackDateTime
    $ws.Range("K3").Value = $HandbackDateTime

    # Rebuild hyperlinks so the new I2/I3 links land right after their A2/A3
    # siblings (same relative order Excel produces when the links are added
    # interleaved row by row), then restore the shared "HyperLink" look.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98cf1e538dbf1d988282a9d2bd660cd0ce92c087/e2e/b9a358c7-8e1c-46ea-94ea-a5e869538c3e.md", "", "", "b9a358c7-8e1c-46ea-94ea-a5e869538c3e.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98cf1e538dbf1d988282a9d2bd660cd0ce92c087/e2e/b9a358c7-8e1c-46ea-94ea-a5e869538c3e.md", "", "", "b9a358c7-8e1c-46ea-94ea-a5e869538c3e.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98cf1e538dbf1d988282a9d2bd660cd0ce92c087/e2e/e0a40405-df7f-4f03-93b3-560f8c411baf.md", "", "", "e0a40405-df7f-4f03-93b3-560f8c411baf.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98cf1e538dbf1d988282a9d2bd660cd0ce92c087/e2e/e0a40405-df7f-4f03-93b3-560f8c411baf.md", "", "", "e0a40405-df7f-4f03-93b3-560f8c411baf.md")

    $ws.Range("A2").Style = "HyperLink"
    $ws.Range("A3").Style = "HyperLink"
    $ws.Range("I2").Style = "HyperLink"
    $ws.Range("I3").Style = "HyperLink"

    # Widen the columns that now hold the long generated file names.
    $ws.Range("C1").ColumnWidth = 29.9777047293527
    $ws.Range("I1").ColumnWidth = 40
    $ws.Range("J1").ColumnWidth = 40
}

Update-LanguageSheet "zh-cn" `
    "b9a358c7-8e1c-46ea-94ea-a5e869538c3e.6f9af68845285ad22dd482a18351f42a981dbe13.zh-cn.xlf" `
    "e0a40405-df7f-4f03-93b3-560f8c411baf.b2df989fc2c2d61e9271556d0df96963afb0e84a.zh-cn.xlf" `
    "2016-08-30 00:31:15"

Update-LanguageSheet "de-de" `
    "b9a358c7-8e1c-46ea-94ea-a5e869538c3e.6f9af68845285ad22dd482a18351f42a981dbe13.de-de.xlf" `
    "e0a40405-df7f-4f03-93b3-560f8c411baf.b2df989fc2c2d61e9271556d0df96963afb0e84a.de-de.xlf" `
    "2016-08-30 00:31:22"
